$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36; this shifts the existing rows 36 and
# 37 down to 37 and 38 respectively.
$ws.Rows.Item(36).Insert()

# Pull the border/style formatting (style index 1 -> thin box border) from
# row 35 (a fully-populated A:F row) down into the new row 36 so every one
# of its six cells picks up the same look the rest of the table uses.
$ws.Range("A35:F35").Copy()
$ws.Range("A36:F36").PasteSpecial(-4122)

# Plain text values (not numeric-looking) can be assigned directly - Excel
# keeps these as text automatically.
$ws.Range("B36").Value = "APEL CHERRY/STRAWBRY"
$ws.Range("C36").Value = "RBI01S"
$ws.Range("F36").Value = "PT,(E-1H)"

# The remaining new values look like plain numbers ("20129832", "2", "14"),
# so a direct .Value assignment on a General-formatted cell would store them
# as numbers instead of text. Route them through a text formula first, then
# collapse the formula down to its literal value with a values-only paste -
# this keeps the cell's existing style/number-format untouched while still
# landing a text (shared-string) value in the cell, matching how the rest of
# this column is stored.
$ws.Range("A36").Formula = '="20129832"'
$ws.Range("A36").Copy()
$ws.Range("A36").PasteSpecial(-4163)

$ws.Range("D36").Formula = '="2"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)

$ws.Range("E36").Formula = '="14"'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
